# Apply the updated crypto price/volume snapshot (GitHub Actions refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D ("Price") values that read as plain numbers need a leading
# apostrophe (just like typing '1.003 into Excel) so they stay text,
# matching the original inline-string cells instead of becoming numbers.

$ws.Range("D2").Value = '25.989.41'
$ws.Range("E2").Value = '  +0.12%  '
$ws.Range("D3").Value = '1.640.03'
$ws.Range("E3").Value = '  -0.61%  '
$ws.Range("D4").Value = '''1.003'
$ws.Range("E4").Value = '  -0.45%  '
$ws.Range("D5").Value = '''214.88'
$ws.Range("E5").Value = '  -0.63%  '
$ws.Range("D6").Value = '''0.5093'
$ws.Range("E6").Value = '  -0.55%  '
$ws.Range("D7").Value = '''1.002'
$ws.Range("E7").Value = '  -0.55%  '
$ws.Range("E8").Value = '  -0.25%  '
$ws.Range("D9").Value = '''0.06353'
$ws.Range("E9").Value = '  -1.34%  '
$ws.Range("D10").Value = '''19.84'
$ws.Range("E10").Value = '  +0.50%  '
$ws.Range("D11").Value = '''0.07760'
$ws.Range("E11").Value = '  -0.32%  '
$ws.Range("D12").Value = '''4.272'
$ws.Range("E12").Value = '  -1.38%  '
$ws.Range("D13").Value = '1.634.69'
$ws.Range("E13").Value = '  -0.65%  '
$ws.Range("D14").Value = '''0.5466'
$ws.Range("E14").Value = '  -0.42%  '
$ws.Range("D15").Value = '0.0₅7747'
$ws.Range("E15").Value = '  -2.08%  '
$ws.Range("D16").Value = '''64.32'
$ws.Range("E16").Value = '  -1.10%  '
$ws.Range("D17").Value = '26.004.05'
$ws.Range("E17").Value = '  -0.16%  '
$ws.Range("D18").Value = '''1.002'
$ws.Range("E18").Value = '  -0.51%  '
$ws.Range("D19").Value = '''195.80'
$ws.Range("E19").Value = '  -1.83%  '
$ws.Range("D20").Value = '''4.429'
$ws.Range("E20").Value = '  -0.88%  '
$ws.Range("D21").Value = '''9.917'
$ws.Range("E21").Value = '  -1.48%  '
$ws.Range("D22").Value = '''6.085'
$ws.Range("E22").Value = '  +0.09%  '
$ws.Range("E23").Value = '  -0.70%  '
$ws.Range("D24").Value = '''1.891'
$ws.Range("E24").Value = '  +1.16%  '
$ws.Range("D25").Value = '''143.19'
$ws.Range("E25").Value = '  +1.91%  '
$ws.Range("D26").Value = '''0.1232'
$ws.Range("E26").Value = '  +6.88%  '
$ws.Range("D27").Value = '''6.860'
$ws.Range("E27").Value = '  -0.94%  '
$ws.Range("D28").Value = '''15.57'
$ws.Range("E28").Value = '  -1.52%  '
$ws.Range("E29").Value = '  -0.50%  '
$ws.Range("D30").Value = '''0.04858'
$ws.Range("E30").Value = '  -3.66%  '
$ws.Range("D31").Value = '''3.274'
$ws.Range("E31").Value = '  -0.77%  '
$ws.Range("D32").Value = '''3.220'
$ws.Range("D33").Value = '''1.542'
$ws.Range("E33").Value = '  -0.52%  '
$ws.Range("E34").Value = '  +0.25%  '
$ws.Range("D35").Value = '''0.9140'
$ws.Range("E35").Value = '  +1.78%  '
$ws.Range("D36").Value = '''2.565'
$ws.Range("E36").Value = '  -1.07%  '
$ws.Range("D37").Value = '''0.5552'
$ws.Range("E37").Value = '  -0.35%  '
$ws.Range("D38").Value = '1.095.93'
$ws.Range("E38").Value = '  -3.94%  '
$ws.Range("E39").Value = '  +0.12%  '
$ws.Range("E40").Value = '  -0.68%  '
$ws.Range("D41").Value = '''2.526'
$ws.Range("E41").Value = '  -1.57%  '
$ws.Range("D42").Value = '''5.587'
$ws.Range("E42").Value = '  -1.73%  '
$ws.Range("D43").Value = '''0.8056'
$ws.Range("E43").Value = '  -1.62%  '
$ws.Range("E44").Value = '  -0.80%  '
$ws.Range("D45").Value = '0.0₈120'
$ws.Range("E45").Value = '  -3.64%  '
$ws.Range("D46").Value = '1.779.56'
$ws.Range("E46").Value = '  -0.39%  '
$ws.Range("D47").Value = '''0.4538'
$ws.Range("E47").Value = '  -0.09%  '
$ws.Range("D48").Value = '''1.008'
$ws.Range("E48").Value = '  -0.09%  '
$ws.Range("D49").Value = '''55.23'
$ws.Range("E49").Value = '  -0.44%  '
$ws.Range("D50").Value = '''0.05215'
$ws.Range("E50").Value = '  +2.28%  '
$ws.Range("D51").Value = '''7.524'
$ws.Range("E51").Value = '  +1.37%  '
